# Absenzenlisten-Template 2016/2017 - Semester 2 (Montag)
# "zu grosse Schrift in einzelnen Zellen korrigiert"
#
# A handful of the single-letter "F" (Ferien/holiday) cells in the
# attendance table were left at the document's default run size
# (11 pt) instead of the 10 pt used by every other cell in the table.
# Walk the table, find the cells whose entire content is just "F" and
# whose font size is not yet 10 pt, and correct both the regular and
# complex-script size so they match the rest of the sheet.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$fixedCount = 0

foreach ($cell in $table.Range.Cells) {
    $text = $cell.Range.Text
    # Cell text includes the trailing cell-mark (Chr 7) and, for the
    # last cell in a row, a paragraph mark (Chr 13) before it - strip
    # both so we can compare against the plain letter.
    $plain = $text.TrimEnd([char]7).TrimEnd([char]13)

    if ($plain -eq "F") {
        if ($cell.Range.Font.Size -ne 10) {
            $cell.Range.Font.Size = 10
            $cell.Range.Font.SizeBi = 10
            $fixedCount = $fixedCount + 1
        }
    }
}

Write-Host "Corrected font size on $fixedCount 'F' cell(s)."
